$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last tender row (row 8, "Mise en conformité de la ligne 60 kV ...")
# This shifts nothing else in content - rows 2-7 keep their original objet/date_limite
# values; only the hyperlink URLs in column C need refreshing afterwards.
$ws.Rows.Item(8).Delete()

# Drop all existing hyperlinks (and their relationships) so we can recreate them
# with the new, per-row target URLs.
$ws.Hyperlinks.Delete()

# New distinct consultation URLs for each remaining row.
$urls = @{
    2 = "https://achats.adm.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=4410&orgAcronyme=e3r&echanges"
    3 = "https://achats.adm.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=4409&orgAcronyme=e3r&echanges"
    4 = "https://achats.adm.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=4387&orgAcronyme=e3r&echanges"
    5 = "https://achats.adm.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=4408&orgAcronyme=e3r&echanges"
    6 = "https://achats.adm.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=4403&orgAcronyme=e3r&echanges"
    7 = "https://achats.adm.co.ma/?page=entreprise.EntrepriseDetailsConsultation&refConsultation=4386&orgAcronyme=e3r&echanges"
}

foreach ($row in 2..7) {
    $url = $urls[$row]
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $url
    $ws.Hyperlinks.Add($cell, $url)
    $cell.Style = "Hyperlink"
}
